$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the report title (report date 2025-09-02 -> 2025-09-03)
$ws.Range("A1").Value = "萊爾富 工作統計表  篩選月份：202509   (  製表日期:2025-09-03  )"

# Turn word-wrap back on for P25 / AC25 (they'd lost it, unlike every
# other data row's "工作內容/報修說明" columns which wrap)
$ws.Range("P25").WrapText = $true
$ws.Range("AC25").WrapText = $true

# Duplicate the formatting of the row above onto the new row 26
$ws.Range("A24:AK24").Copy()
$ws.Range("A26:AK26").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new maintenance-report row
$ws.Range("A26").Value = 24
$ws.Range("B26").Value = "服務"
$ws.Range("C26").Value = 2025090599
$ws.Range("F26").Value = 4658
$ws.Range("G26").Value = "淡水英才店"
$ws.Range("H26").Value = "新北市淡水區"
$ws.Range("Q26").Value = "THILF04658"
$ws.Range("R26").Value = "新北一"
$ws.Range("S26").Value = "吳宗鴻"
$ws.Range("T26").Value = 1
$ws.Range("U26").Value = "已完工"
$ws.Range("V26").Value = "2025-09-03 12:48:12"
$ws.Range("W26").Value = "2025-09-03 12:20:00"
$ws.Range("X26").Value = "2025-09-03 12:47:00"
$ws.Range("Z26").Value = 0.5
$ws.Range("AB26").Value = "到場處理"
$ws.Range("AC26").Value = "PMQ3+TVV"
$ws.Range("AD26").Value = "O"
$ws.Range("AK26").Value = "O"

# Update print area to include the new row
$ws.PageSetup.PrintArea = '$A$1:$AK$26'

# Move the active selection (matches the author's final cursor position)
$ws.Range("AC23").Select()
